$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update corrected values for rows 2-6

# Row 2
$ws.Range("D2").Value = 975
$ws.Range("E2").Value = -24
$ws.Range("F2").Value = -24
$ws.Range("G2").Value = -82
$ws.Range("H2").Value = -70
$ws.Range("I2").Value = -69
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 910
$ws.Range("L2").Value = 684
$ws.Range("M2").Value = 226
$ws.Range("N2").Value = 226
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 51
$ws.Range("Q2").Value = -2
$ws.Range("R2").Value = 6
$ws.Range("S2").Value = -3
$ws.Range("T2").Value = 2
$ws.Range("U2").Value = -4
$ws.Range("V2").Value = 327
$ws.Range("W2").Value = -2.47
$ws.Range("X2").Value = -7.18
$ws.Range("Y2").Value = -26.58
$ws.Range("Z2").Value = -7.39
$ws.Range("AA2").Value = 302.32
$ws.Range("AB2").Value = 331.65
$ws.Range("AC2").Value = -560
$ws.Range("AD2").Value = -1.87
$ws.Range("AE2").Value = 1850
$ws.Range("AF2").Value = 0.5600000000000001
$ws.Range("AG2").Value = 4
$ws.Range("AH2").Value = 0.4
$ws.Range("AI2").Value = -0.73
$ws.Range("AJ2").Value = 12407093

# Row 3
$ws.Range("D3").Value = 916
$ws.Range("E3").Value = -25
$ws.Range("F3").Value = -25
$ws.Range("G3").Value = -45
$ws.Range("H3").Value = -40
$ws.Range("I3").Value = -40
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 831
$ws.Range("L3").Value = 647
$ws.Range("M3").Value = 184
$ws.Range("N3").Value = 184
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 51
$ws.Range("Q3").Value = -8
$ws.Range("R3").Value = 31
$ws.Range("S3").Value = -24
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = -9
$ws.Range("V3").Value = 304
$ws.Range("W3").Value = -2.69
$ws.Range("X3").Value = -4.41
$ws.Range("Y3").Value = -19.61
$ws.Range("Z3").Value = -4.65
$ws.Range("AA3").Value = 351.32
$ws.Range("AB3").Value = 251.63
$ws.Range("AC3").Value = -324
$ws.Range("AD3").Value = -4.98
$ws.Range("AE3").Value = 1507
$ws.Range("AF3").Value = 1.07
$ws.Range("AG3").Value = 4
$ws.Range("AH3").Value = 0.26
$ws.Range("AI3").Value = -1.25
$ws.Range("AJ3").Value = 12407093

# Row 4
$ws.Range("D4").Value = 873
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = -44
$ws.Range("H4").Value = -42
$ws.Range("I4").Value = -43
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 813
$ws.Range("L4").Value = 660
$ws.Range("M4").Value = 153
$ws.Range("N4").Value = 152
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 51
$ws.Range("Q4").Value = 5
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = -8
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 5
$ws.Range("V4").Value = 293
$ws.Range("W4").Value = 0.23
$ws.Range("X4").Value = -4.8
$ws.Range("Y4").Value = -25.32
$ws.Range("Z4").Value = -5.09
$ws.Range("AA4").Value = 432.28
$ws.Range("AB4").Value = 181.52
$ws.Range("AC4").Value = -343
$ws.Range("AD4").Value = -5.08
$ws.Range("AE4").Value = 1225
$ws.Range("AF4").Value = 1.42
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 12407093

# Row 5
$ws.Range("D5").Value = 783
$ws.Range("E5").Value = -39
$ws.Range("F5").Value = -39
$ws.Range("G5").Value = -52
$ws.Range("H5").Value = -29
$ws.Range("I5").Value = -28
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 733
$ws.Range("L5").Value = 554
$ws.Range("M5").Value = 178
$ws.Range("N5").Value = 178
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 84
$ws.Range("Q5").Value = -31
$ws.Range("R5").Value = -5
$ws.Range("S5").Value = 43
$ws.Range("T5").Value = 3
$ws.Range("U5").Value = -35
$ws.Range("V5").Value = 280
$ws.Range("W5").Value = -5.03
$ws.Range("X5").Value = -3.69
$ws.Range("Y5").Value = -17.16
$ws.Range("Z5").Value = -3.74
$ws.Range("AA5").Value = 310.51
$ws.Range("AB5").Value = 106.92
$ws.Range("AC5").Value = -198
$ws.Range("AD5").Value = -4.75
$ws.Range("AE5").Value = 1060
$ws.Range("AF5").Value = 0.89
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 16823040

# Row 6
$ws.Range("D6").Value = 716
$ws.Range("E6").Value = -53
$ws.Range("F6").Value = -53
$ws.Range("G6").Value = -9
$ws.Range("H6").Value = -40
$ws.Range("I6").Value = -38
$ws.Range("K6").Value = 482
$ws.Range("L6").Value = 348
$ws.Range("M6").Value = 134
$ws.Range("N6").Value = 136
$ws.Range("P6").Value = 84
$ws.Range("Q6").Value = -61
$ws.Range("R6").Value = 207
$ws.Range("S6").Value = -152
$ws.Range("T6").Value = 25
$ws.Range("U6").Value = -85
$ws.Range("V6").Value = 128
$ws.Range("W6").Value = -7.4
$ws.Range("X6").Value = -5.59
$ws.Range("Y6").Value = -24.45
$ws.Range("Z6").Value = -6.59
$ws.Range("AA6").Value = 259.29
$ws.Range("AB6").Value = 53.81
$ws.Range("AC6").Value = -228
$ws.Range("AD6").Value = -4.98
$ws.Range("AE6").Value = 806
$ws.Range("AF6").Value = 1.41
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 16823040

# Clear removed values for rows 7-9 (only A/B/C remain)

# Row 7
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
